$wb = $excel.ActiveWorkbook

# --- Remove the now-unused "Parameters" sheet (and its external-reference
#     renumbering is handled automatically by the host on save) ---
[void]$wb.Worksheets.Item("Parameters").Delete()

$ws = $wb.Worksheets.Item(1)

# --- Extend the data table with the 2020-2023 columns (E:H) ---

# Row 3 — year headers
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4 — number of local governments (constant across years).
# D4 used to hold the figure as text (shared string "484"); it becomes a
# genuine number once the series is extended, so reset it here too.
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5 — proportion (%) values
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6 — counts of governments implementing DRR strategies
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- Carry the column-D formatting across into the new columns ---
$ws.Range("D3").Copy()
$ws.Range("E3:H3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("E4:H4").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("E5:H5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6:H6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row-height adjustments for the now-taller wrapped rows ---
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 44.25
$ws.Rows.Item(6).RowHeight = 51.75

# --- Restore the selection left behind in the source workbook ---
[void]$ws.Range("D9").Select()
